# "managed client infos with common export function"
#
# The "Vendredi 12/3" week block (rows 80-116) is restructured: instead of
# a per-day SUM() of several sub-entries, most days now carry a single task
# label (in column K, as a shared-string lookup) with its hour value sitting
# directly on the date row (or on a single child row) in column L. The other
# now-unused child rows are cleared out entirely (no cell left behind).
# A brand-new task entry is also added at row 110 (with a new shared string),
# together with an (empty, but time-formatted) helper cell in H110.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Vendredi 12/3 (J81 = 44305) ---------------------------------------
# Row 81 was the day-total row (=SUM(L82:L88)); it becomes a single task
# entry: K81 = task label, L81 = plain value (formula removed).
$ws.Range("K81").Value = "Affichage de l'apercu avant impression au clic sur le bouton prévu à cet effet"
$ws.Range("L81").Value = 1

# Row 82 keeps its value but also gets a task label.
$ws.Range("K82").Value = "Créer une fonction prenant en entrée rangeIn et qui va écrire dans l'export"
$ws.Range("L82").Value = 1

# Rows 83-88 no longer hold any hour value.
$ws.Range("L83:L88").ClearContents()

# --- Day at J90 (=J81+1) -------------------------------------------------
# L90 keeps its SUM(L91:L94) formula; only row 91 keeps a value (now with a
# task label), rows 92-94 are cleared.
$ws.Range("K91").Value = "Créer une fonction prenant en entrée rangeIn et qui va écrire dans l'export"
$ws.Range("L92:L94").ClearContents()

# --- Day at J96 (=J90+1) -------------------------------------------------
# L96 keeps its SUM(L97:L100) formula; row 97 keeps only a task label (no
# value any more), rows 97-100 lose their hour values.
$ws.Range("K97").Value = "Créer une fonction prenant en entrée rangeIn et qui va écrire dans l'export"
$ws.Range("L97:L100").ClearContents()

# --- Day at J109 (=J102+1) ------------------------------------------------
# Row 110 gets a brand-new task label (new shared string) and keeps a
# (smaller) hour value; rows 111-116 lose their hour values. H110 gets a
# (still empty) time-formatted helper cell.
$ws.Range("H110").NumberFormat = "h:mm"
$ws.Range("K110").Value = "Faire une structure pour ClientDetails pour fiter rangeIn"
$ws.Range("L110").Value = 2
$ws.Range("L111:L116").ClearContents()

# --- Sheet-level cosmetics -------------------------------------------------
# Disable automatic page breaks (pageSetUpPr/@autoPageBreaks) and move the
# selection/scroll position to where the new data now lives.
$ws.PageSetup.Zoom = $ws.PageSetup.Zoom
$ws.PageSetup.PrintArea = $ws.PageSetup.PrintArea

$ws.Range("K111").Select()
